# Refresh market-price-derived columns (H,I,J,K,L,M,N) across the Leve
# profit tables, per the latest scheduled Universalis price pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3705643.8
$ws.Range("J17").Value = 3705643.8
$ws.Range("L17").Value = 11116931.4
$ws.Range("N17").Value = -11117267.4
$ws.Range("H28").Value = 308.07693
$ws.Range("I28").Value = 336.9091
$ws.Range("K28").Value = 336.9091
$ws.Range("M28").Value = 148.0909
$ws.Range("H86").Value = 2784.2632
$ws.Range("I86").Value = 2769.3635
$ws.Range("K86").Value = 2769.3635
$ws.Range("M86").Value = -1646.3635
$ws.Range("H89").Value = 2784.2632
$ws.Range("I89").Value = 2769.3635
$ws.Range("K89").Value = 13846.8175
$ws.Range("M89").Value = -8230.817499999999
$ws.Range("H106").Value = 3277.4443
$ws.Range("I106").Value = 1946.7858
$ws.Range("K106").Value = 1946.7858
$ws.Range("M106").Value = -1315.7858
$ws.Range("H107").Value = 238.61905
$ws.Range("I107").Value = 257.17648
$ws.Range("K107").Value = 257.17648
$ws.Range("M107").Value = 1662.82352
$ws.Range("H111").Value = 600
$ws.Range("I111").Value = 0
$ws.Range("K111").Value = 0
$ws.Range("M111").Value = $null
$ws.Range("H137").Value = 2323.739
$ws.Range("I137").Value = 2275
$ws.Range("K137").Value = 6825
$ws.Range("M137").Value = -4275
$ws.Range("H138").Value = 2687.205
$ws.Range("J138").Value = 3926.95
$ws.Range("L138").Value = 11780.85
$ws.Range("N138").Value = -22060.85

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2126.32
$ws.Range("I2").Value = 1833.0869
$ws.Range("K2").Value = 1833.0869
$ws.Range("M2").Value = -1720.0869
$ws.Range("H32").Value = 5332.4
$ws.Range("I32").Value = 4390.5874
$ws.Range("J32").Value = 34999.5
$ws.Range("K32").Value = 4390.5874
$ws.Range("L32").Value = 34999.5
$ws.Range("M32").Value = -4103.5874
$ws.Range("N32").Value = -35573.5
$ws.Range("H45").Value = 4825.3335
$ws.Range("I45").Value = 3886.077
$ws.Range("K45").Value = 3886.077
$ws.Range("M45").Value = -3509.077
$ws.Range("H61").Value = 10965.111
$ws.Range("I61").Value = 11717
$ws.Range("K61").Value = 11717
$ws.Range("M61").Value = -11505
$ws.Range("H116").Value = 2126.32
$ws.Range("I116").Value = 1833.0869
$ws.Range("K116").Value = 1833.0869
$ws.Range("M116").Value = 460.9131
$ws.Range("H122").Value = 1828.75
$ws.Range("I122").Value = 1805
$ws.Range("K122").Value = 5415
$ws.Range("M122").Value = -2965
$ws.Range("H132").Value = 1587.8276
$ws.Range("I132").Value = 1448.5186
$ws.Range("K132").Value = 4345.5558
$ws.Range("M132").Value = -1815.5558
$ws.Range("H136").Value = 10965.111
$ws.Range("I136").Value = 11717
$ws.Range("K136").Value = 35151
$ws.Range("M136").Value = -32601

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2126.32
$ws.Range("I3").Value = 1833.0869
$ws.Range("K3").Value = 1833.0869
$ws.Range("M3").Value = -1719.0869
$ws.Range("H86").Value = 1618.4286
$ws.Range("I86").Value = 1629.5555
$ws.Range("K86").Value = 1629.5555
$ws.Range("M86").Value = -506.5554999999999
$ws.Range("H89").Value = 1618.4286
$ws.Range("I89").Value = 1629.5555
$ws.Range("K89").Value = 8147.7775
$ws.Range("M89").Value = -2531.7775
$ws.Range("H105").Value = 1834.5333
$ws.Range("I105").Value = 2254.111
$ws.Range("K105").Value = 2254.111
$ws.Range("M105").Value = -507.1109999999999
$ws.Range("H107").Value = 1586.5
$ws.Range("I107").Value = 1329.25
$ws.Range("J107").Value = 2872.75
$ws.Range("K107").Value = 1329.25
$ws.Range("L107").Value = 2872.75
$ws.Range("M107").Value = 590.75
$ws.Range("N107").Value = -6712.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").Value = $null
$ws.Range("H52").Value = 93084.5
$ws.Range("J52").Value = 93084.5
$ws.Range("L52").Value = 93084.5
$ws.Range("N52").Value = -93672.5
$ws.Range("H58").Value = 11408.6875
$ws.Range("I58").Value = 6256.6665
$ws.Range("J58").Value = 14499.9
$ws.Range("K58").Value = 6256.6665
$ws.Range("L58").Value = 14499.9
$ws.Range("M58").Value = -6053.6665
$ws.Range("N58").Value = -14905.9
$ws.Range("H107").Value = 1542.1111
$ws.Range("I107").Value = 521.5714
$ws.Range("K107").Value = 521.5714
$ws.Range("M107").Value = 1398.4286
$ws.Range("H132").Value = 5655.727
$ws.Range("I132").Value = 3449.1428
$ws.Range("J132").Value = 9517.25
$ws.Range("K132").Value = 10347.4284
$ws.Range("L132").Value = 28551.75
$ws.Range("M132").Value = -7817.428400000001
$ws.Range("N132").Value = -33611.75
$ws.Range("H134").Value = 2606.7407
$ws.Range("I134").Value = 1666.2858
$ws.Range("K134").Value = 4998.857400000001
$ws.Range("M134").Value = -2463.857400000001
$ws.Range("H136").Value = 11408.6875
$ws.Range("I136").Value = 6256.6665
$ws.Range("J136").Value = 14499.9
$ws.Range("K136").Value = 18769.9995
$ws.Range("L136").Value = 43499.7
$ws.Range("M136").Value = -16219.9995
$ws.Range("N136").Value = -48599.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 850
$ws.Range("I22").Value = 500
$ws.Range("J22").Value = 1200
$ws.Range("K22").Value = 1500
$ws.Range("L22").Value = 3600
$ws.Range("M22").Value = -1331
$ws.Range("N22").Value = -3938
$ws.Range("H27").Value = 850
$ws.Range("I27").Value = 500
$ws.Range("J27").Value = 1200
$ws.Range("K27").Value = 1500
$ws.Range("L27").Value = 3600
$ws.Range("M27").Value = -1398
$ws.Range("N27").Value = -3804
$ws.Range("H81").Value = 8942.714
$ws.Range("I81").Value = 1649.75
$ws.Range("K81").Value = 4949.25
$ws.Range("M81").Value = -3826.25
$ws.Range("H84").Value = 8942.714
$ws.Range("I84").Value = 1649.75
$ws.Range("K84").Value = 14847.75
$ws.Range("M84").Value = -9231.75
$ws.Range("H88").Value = 15000
$ws.Range("J88").Value = 15000
$ws.Range("L88").Value = 45000
$ws.Range("N88").Value = -45856
$ws.Range("H91").Value = 15000
$ws.Range("J91").Value = 15000
$ws.Range("L91").Value = 45000
$ws.Range("N91").Value = -47964
$ws.Range("H129").Value = 740709.6
$ws.Range("J129").Value = 1442785.2
$ws.Range("L129").Value = 4328355.6
$ws.Range("N129").Value = -4338355.6
$ws.Range("H134").Value = 1724.75
$ws.Range("I134").Value = 1724.75
$ws.Range("K134").Value = 5174.25
$ws.Range("M134").Value = -104.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10130
$ws.Range("J70").Value = 14798.125
$ws.Range("L70").Value = 14798.125
$ws.Range("N70").Value = -15338.125
$ws.Range("H73").Value = 10130
$ws.Range("J73").Value = 14798.125
$ws.Range("L73").Value = 14798.125
$ws.Range("N73").Value = -16670.125
$ws.Range("H97").Value = 1008.65216
$ws.Range("I97").Value = 822.1111
$ws.Range("K97").Value = 822.1111
$ws.Range("M97").Value = -326.1111
$ws.Range("H102").Value = 3820.087
$ws.Range("I102").Value = 2519.1052
$ws.Range("K102").Value = 2519.1052
$ws.Range("M102").Value = -897.1052
$ws.Range("H109").Value = 69666
$ws.Range("J109").Value = 69666
$ws.Range("L109").Value = 69666
$ws.Range("N109").Value = -71746
$ws.Range("H132").Value = 5933.64
$ws.Range("I132").Value = 5578.9565
$ws.Range("K132").Value = 16736.8695
$ws.Range("M132").Value = -14206.8695

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 373.27274
$ws.Range("I16").Value = 335.6
$ws.Range("K16").Value = 335.6
$ws.Range("M16").Value = -165.6
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").Value = $null
$ws.Range("H100").Value = 4179.5713
$ws.Range("I100").Value = 3949.3333
$ws.Range("J100").Value = 4594
$ws.Range("K100").Value = 3949.3333
$ws.Range("L100").Value = 4594
$ws.Range("M100").Value = -3408.3333
$ws.Range("N100").Value = -5676
$ws.Range("H136").Value = 29855660
$ws.Range("I136").Value = 16398516
$ws.Range("K136").Value = 49195548
$ws.Range("M136").Value = -49192998

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 1050
$ws.Range("I23").Value = 1050
$ws.Range("K23").Value = 1050
$ws.Range("M23").Value = -821
$ws.Range("H100").Value = 607.9091
$ws.Range("J100").Value = 788
$ws.Range("L100").Value = 1576
$ws.Range("N100").Value = -2658
$ws.Range("H132").Value = 7437.5483
$ws.Range("I132").Value = 6337.231
$ws.Range("K132").Value = 19011.693
$ws.Range("M132").Value = -16481.693
$ws.Range("H136").Value = 5554.636
$ws.Range("I136").Value = 3455.7778
$ws.Range("K136").Value = 10367.3334
$ws.Range("M136").Value = -7817.3334
